# Card10 - add new service event (row 24), 2025-12-08, by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

$bs = [char]0x5C

# Previously-blank cells on row 23 were backfilled with the literal text "nan"
# (matches how every other fully-populated row in this sheet marks "no data").
$ws.Range("B23").Value = "nan"
$ws.Range("C23").Value = "nan"
$ws.Range("D23").Value = "nan"
$ws.Range("E23").Value = "nan"
$ws.Range("F23").Value = "nan"
$ws.Range("G23").Value = "nan"
$ws.Range("H23").Value = "nan"
$ws.Range("I23").Value = "nan"
$ws.Range("J23").Value = "nan"
$ws.Range("K23").Value = "nan"
$ws.Range("M23").Value = "nan"

# New row 24: new service/maintenance event entry for the machine (card = 10)
# "10" is stored as text (matches column A formatting for every other row),
# so it is entered with a leading apostrophe to force text rather than a number.
$ws.Range("A24").Value = "'10"

$ws.Range("L24").Value = "10${bs}7${bs}2025"
$ws.Range("M24").Value = "9377 h  736 t"
$ws.Range("N24").Value = "تم تغييرزيت الجيربوكس"
$ws.Range("O24").Value = "تيم العمل"
